$d = $word.ActiveDocument

$replacements = @(
    @("71+6=", "16+70="),
    @("20+29=", "84-22="),
    @("27-24=", "16+14="),
    @("93-91=", "87-29="),
    @("32+15=", "38+22="),
    @("3+4=", "74-31="),
    @("57-40=", "17-9="),
    @("32-29=", "41-20="),
    @("13+79=", "82-20="),
    @("81-12=", "87+9="),
    @("68-13=", "26+39="),
    @("79-13=", "10+64="),
    @("55-24=", "31-22="),
    @("95-41=", "63+3="),
    @("1+49=", "47+41="),
    @("99-40=", "3+9="),
    @("74-6=", "66-35="),
    @("17+80=", "87-48="),
    @("27+52=", "59-1="),
    @("0+97=", "57-21="),
    @("11+46=", "77-59="),
    @("18+9=", "2+36="),
    @("87-55=", "94-62="),
    @("3+19=", "13+1="),
    @("51-27=", "44-15="),
    @("6+61=", "30-5="),
    @("24+6=", "80+11="),
    @("97-49=", "20+30="),
    @("39-39=", "1+43="),
    @("99-2=", "58-49="),
    @("11+45=", "89-87="),
    @("33-22=", "76-4="),
    @("48+16=", "87-21="),
    @("96-70=", "59-8="),
    @("95-70=", "2+20="),
    @("43+9=", "4-0="),
    @("9+36=", "75+23="),
    @("47+17=", "11+53="),
    @("38+42=", "91-69="),
    @("91-78=", "42-17="),
    @("81-36=", "50-37="),
    @("22+43=", "10+71="),
    @("92-41=", "61+37="),
    @("49-27=", "43+23="),
    @("49+45=", "20+54="),
    @("36+62=", "74+17="),
    @("29+40=", "48-13="),
    @("85-50=", "76-7="),
    @("44+50=", "97-51="),
    @("52+33=", "51-3="),
    @("78-49=", "92-74="),
    @("18+0=", "50-26="),
    @("53-8=", "62-41="),
    @("27+20=", "77-68="),
    @("32+7=", "55-51="),
    @("39+27=", "3+35="),
    @("53+23=", "77-1="),
    @("91-65=", "19-6="),
    @("43+5=", "71-66="),
    @("66+8=", "47-13="),
    @("61-38=", "51+17="),
    @("81-31=", "33+13="),
    @("74+12=", "56+17="),
    @("77-70=", "15+74="),
    @("26+58=", "14+29="),
    @("43+33=", "31+14="),
    @("40+14=", "69+27="),
    @("1+84=", "38+12="),
    @("29-17=", "74+9="),
    @("37-3=", "33+46="),
    @("10+14=", "42+16="),
    @("54-4=", "31-3="),
    @("39-29=", "77-71="),
    @("41+52=", "51-21="),
    @("35-6=", "89-44="),
    @("31+34=", "16+83="),
    @("67-23=", "53-29="),
    @("77+19=", "97-9="),
    @("5+66=", "36+31="),
    @("74-9=", "86-9="),
    @("59-35=", "32+51="),
    @("67-25=", "41+27="),
    @("57-51=", "73-9="),
    @("24+25=", "48+6="),
    @("49+19=", "58+12="),
    @("28+60=", "45+9="),
    @("16+24=", "45-37="),
    @("86-39=", "42+13="),
    @("56-47=", "28+65="),
    @("60+32=", "31+32="),
    @("1+54=", "15-13="),
    @("80-26=", "92-89="),
    @("47+11=", "64-36="),
    @("13-5=", "73-33="),
    @("92+7=", "98-95="),
    @("3+50=", "94-50="),
    @("49-45=", "86-29="),
    @("91-71=", "31+38="),
    @("25-17=", "55+17="),
    @("5+3=", "19+33="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Done: applied $($replacements.Count) replacements"
